# TC29_Canine_Filter_Breed-Gordon.xlsx — "startup" sheet fix-up
#
# The CasesTab query (row 2 / cell B2) had an extra, erroneous `Cohort`
# projection tacked on to the end of the Cypher query. Drop that trailing
# line (and the now-dangling comma before it) so the query matches the
# other corrected TC* queries in this batch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN['Gordon Setter']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# Row heights settle back down now that the CasesTab query lost a couple of
# wrapped lines (SamplesTab/FilesTab text is unchanged, but widths/heights
# were re-measured on the machine that saved this revision).
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.33333333333333
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# Selection moves from B4 (last saved cursor position) to B2.
$ws.Range("B2").Select()
